# Updates cryptocurrency price/volume(1h) figures (and, for rows 46/47,
# swaps the Monero / InjectiveProtocol entries) to match the latest
# GitHub Actions scrape of coinranking.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All the cells we touch hold plain text (prices/percentages are
# formatted strings, not numbers) in the source workbook. Assigning
# a numeric-looking string straight to .Value lets Excel's COM layer
# silently reinterpret it as a number (e.g. "1.00" -> 1, dropping the
# trailing zero; or mangling the subscripted PEPE price). Routing the
# write through a quoted formula and then collapsing it back down to
# a literal with a values-only paste keeps every value as text, byte
# for byte, with no left-over formula or number formatting behind it.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Formula = "=""" + $text + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# Row 2
Set-TextValue "D2" "63.080.17"
Set-TextValue "E2" "  +0.53%  "

# Row 3
Set-TextValue "D3" "3.070.35"
Set-TextValue "E3" "  +0.04%  "

# Row 4
Set-TextValue "E4" "  +0.12%  "

# Row 5
Set-TextValue "D5" "537.94"
Set-TextValue "E5" "  -0.59%  "

# Row 6
Set-TextValue "D6" "137.18"
Set-TextValue "E6" "  +2.52%  "

# Row 7
Set-TextValue "E7" "  -0.04%  "

# Row 8
Set-TextValue "D8" "3.063.91"
Set-TextValue "E8" "  +0.09%  "

# Row 9
Set-TextValue "D9" "0.492"
Set-TextValue "E9" "  +0.64%  "

# Row 10
Set-TextValue "E10" "  +0.68%  "

# Row 11
Set-TextValue "D11" "6.21"
Set-TextValue "E11" "  +0.76%  "

# Row 12
Set-TextValue "E12" "  -2.12%  "

# Row 13
Set-TextValue "E13" "  +0.45%  "

# Row 14
Set-TextValue "D14" "34.42"
Set-TextValue "E14" "  -1.22%  "

# Row 15
Set-TextValue "D15" "3.562.59"
Set-TextValue "E15" "  +1.62%  "

# Row 16
Set-TextValue "D16" "63.080.99"
Set-TextValue "E16" "  +0.74%  "

# Row 17
Set-TextValue "D17" "0.112"
Set-TextValue "E17" "  +1.53%  "

# Row 18
Set-TextValue "D18" "3.068.18"
Set-TextValue "E18" "  +0.27%  "

# Row 19
Set-TextValue "E19" "  -0.88%  "

# Row 20
Set-TextValue "D20" "469.74"
Set-TextValue "E20" "  -2.18%  "

# Row 21
Set-TextValue "D21" "13.46"
Set-TextValue "E21" "  +0.16%  "

# Row 22
Set-TextValue "D22" "0.694"
Set-TextValue "E22" "  -2.30%  "

# Row 23
Set-TextValue "D23" "7.01"
Set-TextValue "E23" "  -2.80%  "

# Row 24
Set-TextValue "D24" "78.41"
Set-TextValue "E24" "  -0.52%  "

# Row 25
Set-TextValue "D25" "12.13"
Set-TextValue "E25" "  +0.63%  "

# Row 26
Set-TextValue "E26" "  +0.21%  "

# Row 27
Set-TextValue "E27" "  -0.81%  "

# Row 28
Set-TextValue "D28" "7.87"
Set-TextValue "E28" "  -4.89%  "

# Row 29
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.26%  "

# Row 30
Set-TextValue "D30" "26.11"
Set-TextValue "E30" "  -0.34%  "

# Row 31
Set-TextValue "D31" "1.15"
Set-TextValue "E31" "  +4.83%  "

# Row 32
Set-TextValue "E32" "  -2.96%  "

# Row 33
Set-TextValue "D33" "59.02"
Set-TextValue "E33" "  +1.79%  "

# Row 34
Set-TextValue "E34" "  -5.21%  "

# Row 35
Set-TextValue "D35" "5.47"
Set-TextValue "E35" "  +6.26%  "

# Row 36
Set-TextValue "D36" "5.94"
Set-TextValue "E36" "  -0.69%  "

# Row 37
Set-TextValue "D37" "481.73"
Set-TextValue "E37" "  -2.01%  "

# Row 38
Set-TextValue "D38" "3.260.94"
Set-TextValue "E38" "  +3.73%  "

# Row 39
Set-TextValue "D39" "0.0396"
Set-TextValue "E39" "  +0.63%  "

# Row 40
Set-TextValue "D40" "0.0793"
Set-TextValue "E40" "  -0.94%  "

# Row 41
Set-TextValue "E41" "  +0.43%  "

# Row 42
Set-TextValue "E42" "  +0.51%  "

# Row 43
Set-TextValue "E43" "  +1.17%  "

# Row 44
Set-TextValue "E44" "  -1.39%  "

# Row 45
Set-TextValue "E45" "  +0.11%  "

# Row 46
Set-TextValue "B46" "InjectiveProtocol"
Set-TextValue "C46" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D46" "25.32"
Set-TextValue "E46" "  +2.25%  "

# Row 47
Set-TextValue "B47" "Monero"
Set-TextValue "C47" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D47" "123.37"
Set-TextValue "E47" "  +4.00%  "

# Row 48
Set-TextValue "E48" "  -2.11%  "

# Row 49
Set-TextValue "D49" "0.108"
Set-TextValue "E49" "  +0.75%  "

# Row 50
Set-TextValue "D50" "0.0₃0520"
Set-TextValue "E50" "  +2.10%  "

# Row 51
Set-TextValue "E51" "  -0.47%  "
